$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $helper = $ws.Range("Z1")
    $escaped = $val.Replace('"', '""')
    $helper.Formula = '="' + $escaped + '"'
    $helper.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $helper.Clear()
}

Set-TextValue "D2" '60.609.98'
$ws.Range("E2").Value = '  -2.16%  '
Set-TextValue "D3" '2.326.22'
$ws.Range("E3").Value = '  -5.40%  '
$ws.Range("E4").Value = '  +0.63%  '
Set-TextValue "D5" '540.74'
Set-TextValue "D6" '135.17'
$ws.Range("E6").Value = '  -7.80%  '
Set-TextValue "D7" '0.994'
$ws.Range("E7").Value = '  -0.70%  '
$ws.Range("E8").Value = '  -11.09%  '
$ws.Range("E9").Value = '  -4.82%  '
$ws.Range("E10").Value = '  -2.63%  '
$ws.Range("E11").Value = '  +0.02%  '
$ws.Range("E12").Value = '  -2.74%  '
$ws.Range("E13").Value = '  -3.20%  '
$ws.Range("E14").Value = '  -6.57%  '
$ws.Range("E15").Value = '  -5.15%  '
$ws.Range("E16").Value = '  -5.65%  '
Set-TextValue "D17" '50.394.38'
$ws.Range("E17").Value = '  -18.49%  '
$ws.Range("E18").Value = '  -9.68%  '
Set-TextValue "D19" '10.51'
$ws.Range("E19").Value = '  -3.73%  '
$ws.Range("E20").Value = '  -2.15%  '
Set-TextValue "D21" '313.41'
$ws.Range("E21").Value = '  -2.14%  '
$ws.Range("E22").Value = '  -6.98%  '
$ws.Range("E23").Value = '  -0.23%  '
Set-TextValue "D24" '62.67'
$ws.Range("E24").Value = '  -2.01%  '
$ws.Range("E25").Value = '  -10.22%  '
Set-TextValue "D26" '8.30'
$ws.Range("E26").Value = '  +5.85%  '
$ws.Range("E27").Value = '  -0.98%  '
$ws.Range("E28").Value = '  -5.27%  '
$ws.Range("E29").Value = '  -5.03%  '
$ws.Range("E30").Value = '  -6.48%  '
Set-TextValue "D31" '1.36'
$ws.Range("E31").Value = '  -8.07%  '
$ws.Range("E32").Value = '  -11.66%  '
Set-TextValue "D33" '0.143'
$ws.Range("E33").Value = '  -2.19%  '
Set-TextValue "D34" '1.77'
$ws.Range("E34").Value = '  -6.37%  '
$ws.Range("E35").Value = '  -7.64%  '
$ws.Range("E36").Value = '  -0.35%  '
Set-TextValue "D37" '4.56'
$ws.Range("E37").Value = '  -4.47%  '
Set-TextValue "D38" '18.40'
$ws.Range("E38").Value = '  +1.00%  '
$ws.Range("E39").Value = '  -3.22%  '
Set-TextValue "D40" '5.16'
$ws.Range("E40").Value = '  -9.49%  '
$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D41" '141.13'
$ws.Range("E41").Value = '  +0.93%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D42" '1.74'
$ws.Range("E42").Value = '  -1.13%  '
$ws.Range("E43").Value = '  -0.07%  '
Set-TextValue "D44" '39.81'
$ws.Range("E44").Value = '  -1.41%  '
Set-TextValue "D45" '139.48'
$ws.Range("E45").Value = '  -2.78%  '
$ws.Range("E46").Value = '  -2.52%  '
$ws.Range("E47").Value = '  -10.77%  '
$ws.Range("E48").Value = '  -3.83%  '
Set-TextValue "D49" '19.01'
$ws.Range("E49").Value = '  -11.62%  '
Set-TextValue "D50" '0.563'
$ws.Range("E50").Value = '  -4.75%  '
Set-TextValue "D51" '0.0890'
$ws.Range("E51").Value = '  -4.66%  '
